$wb = $excel.ActiveWorkbook

# Map of worksheet name -> { row -> new "F" (want-to-go count) value }
# These values reflect a refreshed data scrape (gh-pages output) as of commit 456a3b4.
$updates = @{}
$updates["展览"] = @{
    2 = 192
    3 = 404
    4 = 1138
    5 = 40
    7 = 24
    8 = 1062
    10 = 339
    11 = 419
    13 = 311
    14 = 352
    16 = 60
    17 = 459
    18 = 442
    19 = 5577
    21 = 1557
    22 = 365
    23 = 4741
    25 = 83
    26 = 1490
    27 = 13
    30 = 59
}
$updates["演出"] = @{
    4 = 9
    5 = 130
    8 = 96
    15 = 84
}
$updates["本地生活"] = @{
    2 = 9381
    4 = 2125
}
$updates["全部类型"] = @{
    2 = 9381
    4 = 2125
    5 = 192
    6 = 404
    7 = 1138
    8 = 40
    10 = 24
    11 = 1062
    12 = 339
    13 = 419
    15 = 311
    16 = 352
    18 = 60
    22 = 442
    23 = 5577
    25 = 1557
    28 = 365
    31 = 4741
    33 = 83
    34 = 1490
    35 = 13
    38 = 59
    41 = 84
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowMap = $updates[$sheetName]
    foreach ($row in $rowMap.Keys) {
        $ws.Cells.Item($row, 6).Value = $rowMap[$row]
    }
}
